# jcoleman_Backlog.docx — "Week 6" section rewrite.
#
# The three old Week-6 bullets:
#   "Tweak up the current trees and add more"
#   "Add knock knock joke"
#   "Get started on the poem converter "  (carries the _GoBack bookmark)
#
# become six new bullets (two of them newly highlighted "spell-checked"
# bot names, a new "Week 7" heading, and a rewritten last line that keeps
# the _GoBack bookmark):
#   "Find a jokebot"                (green highlight, "jokebot" flagged)
#   "Learn about the chatbot"       (yellow highlight, "chatbot" flagged)
#   "Get started on the pun engine" (green highlight)
#   "Week 7"                        (underlined heading)
#   "Finish the pun engine"
#   "….Find a new objective off of the proposal"  (keeps _GoBack bookmark)

$d = $word.ActiveDocument
$wdParagraph = 4

# --- Locate the span to replace by content, not by hard-coded index ---
$startRng = $d.Content
$null = $startRng.Find.Execute("Tweak up the current trees and add more", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$startRng.Expand($wdParagraph) | Out-Null

$endRng = $d.Content
$null = $endRng.Find.Execute("Get started on the poem converter", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endRng.Expand($wdParagraph) | Out-Null

$target = $d.Range($startRng.Start, $endRng.End)

# --- Build the replacement WordprocessingML (6 paragraphs) ---
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$p1 = "<w:p $wns>" + `
        "<w:r><w:rPr><w:highlight w:val=`"green`"/></w:rPr>" + `
        "<w:t xml:space=`"preserve`">Find a </w:t></w:r>" + `
        "<w:proofErr w:type=`"spellStart`"/>" + `
        "<w:r><w:rPr><w:highlight w:val=`"green`"/></w:rPr>" + `
        "<w:t>jokebot</w:t></w:r>" + `
        "<w:proofErr w:type=`"spellEnd`"/>" + `
      "</w:p>"

$p2 = "<w:p $wns>" + `
        "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr>" + `
        "<w:t xml:space=`"preserve`">Learn about the </w:t></w:r>" + `
        "<w:proofErr w:type=`"spellStart`"/>" + `
        "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr>" + `
        "<w:t>chatbot</w:t></w:r>" + `
        "<w:proofErr w:type=`"spellEnd`"/>" + `
      "</w:p>"

$p3 = "<w:p $wns>" + `
        "<w:r><w:rPr><w:highlight w:val=`"green`"/></w:rPr>" + `
        "<w:t>Get started on the pun engine</w:t></w:r>" + `
      "</w:p>"

$p4 = "<w:p $wns>" + `
        "<w:pPr><w:rPr><w:u w:val=`"single`"/></w:rPr></w:pPr>" + `
        "<w:r><w:rPr><w:u w:val=`"single`"/></w:rPr>" + `
        "<w:t>Week 7</w:t></w:r>" + `
      "</w:p>"

$p5 = "<w:p $wns>" + `
        "<w:r><w:t>Finish the pun engine</w:t></w:r>" + `
      "</w:p>"

$p6 = "<w:p $wns>" + `
        "<w:r><w:t>&#8230;.Find a new objective off of the proposal</w:t></w:r>" + `
        "<w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/>" + `
        "<w:bookmarkEnd w:id=`"0`"/>" + `
      "</w:p>"

$replacementXml = $p1 + $p2 + $p3 + $p4 + $p5 + $p6

# --- Apply: replacing the whole span in one shot preserves surrounding content ---
[void]$target.InsertXML($replacementXml)
